$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M10").Value = 1521.56
$ws1.Range("M34").Value = "8 de 32"

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F10").Value = 3415.45
$ws2.Range("F34").Value = 34038.58

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 10898.8
$ws3.Range("E16").Value = 10974.3
$ws3.Range("F16").Value = 0.4982741358106533

$ws3.Range("D19").Value = 34556.79
$ws3.Range("E19").Value = -2447.508924442127
$ws3.Range("F19").Value = 1.076224345188009
